$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 1 cell values ---
$ws.Range("A1").Value = 3
$ws.Range("C1").Value = 18
$ws.Range("D1").Value = 14
$ws.Range("E1").Value = 33
$ws.Range("F1").Value = 31
$ws.Range("G1").Value = 27
$ws.Range("H1").Value = 3
$ws.Range("I1").Value = 0.067000000000000004
$ws.Range("J1").Value = 0.021999999999999999
$ws.Range("K1").Value = 0.069999999999999993

# --- Update column widths (columns G:K, i.e. 7-11) ---
# ColumnWidth is expressed in characters and gets quantized by the host
# to the nearest 1/6-character step, so the inputs below are chosen (at
# the midpoint of the relevant quantization bucket, to avoid any
# rounding tie-break ambiguity) as the values that round-trip to the
# closest achievable width to the target raw widths (G=3.140625,
# H=2.140625, I=5.7109375, J=5.7109375, K=4.7109375).
$ws.Columns.Item(7).ColumnWidth = 2.333333333333333
$ws.Columns.Item(8).ColumnWidth = 1.333333333333333
$ws.Columns.Item(9).ColumnWidth = 4.833333333333334
$ws.Columns.Item(10).ColumnWidth = 4.833333333333334
$ws.Columns.Item(11).ColumnWidth = 3.833333333333333
